# Update memorymap with RTC_old
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("memorymap")

# Fill in new rows for mission_time_ss / mission_time_mm / mission_time_hh
$ws.Range("B10").Value = "mission_time_ss"
$ws.Range("E10").Value = "byte"
$ws.Range("F10").Value = ".read()"
$ws.Range("G10").Value = ".update()"

$ws.Range("B11").Value = "mission_time_mm"
$ws.Range("E11").Value = "byte"
$ws.Range("F11").Value = ".read()"
$ws.Range("G11").Value = ".update()"

$ws.Range("B12").Value = "mission_time_hh"
$ws.Range("E12").Value = "byte"
$ws.Range("F12").Value = ".read()"
$ws.Range("G12").Value = ".update()"

# Update the sheet view: remove the frozen topLeftCell / previous selection and set new selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B13").Select()
